$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('10051','1st To 50th Wedding Anniversary Book','2littleboys','2','https://www.notonthehighstreet.com/2littleboys/product/51_the_anniversary_book'),
    @('1014560','Personalised Book Style Novel Notebook','marthabrook','2','https://www.notonthehighstreet.com/marthabrook/product/personalised-book-style-novel-notebook'),
    @('1037522','Graduation Map Card','lisamariedesigns','2','https://www.notonthehighstreet.com/lisamariedesigns/product/graduation-map-card'),
    @('1038914','Personalised Portable Picnic Table Wine Holder','mijmoj','2','https://www.notonthehighstreet.com/mijmoj/product/personalised-picnic-table'),
    @('1047247','Personalised Teacher Thank You Notebook Gift','oliviamorgan','2','https://www.notonthehighstreet.com/oliviamorgan/product/personalised-teacher-thank-you-notebook-gift'),
    @('1062346','Tree Of Life Rose Gold Foil Scarf Gift','studiohop','6','https://www.notonthehighstreet.com/studiohop/product/tree-of-life-foil-scarf-letterbox-gift'),
    @('1067493','Natural Flowers Summer Door Wreath','dibor','3','https://www.notonthehighstreet.com/dibor/product/natural-flowers-autumn-door-wreath'),
    @('1116770','Personalised Christmas Wooden Robin Pair Sign','pinkpineapple','5','https://www.notonthehighstreet.com/pinkpineapple/product/christmas-wooden-robin-pair-sign'),
    @('1150628','Personalised Road Sign','madeforyougifts','6','https://www.notonthehighstreet.com/madeforyougifts/product/personalised-metal-road-sign'),
    @('1162105','','','3',''),
    @('1179458','','','3',''),
    @('1181236','','','2',''),
    @('1182639','','','2',''),
    @('1187498','','','2',''),
    @('1187757','','','2',''),
    @('1198294','','','2',''),
    @('1209526','','','2',''),
    @('1234030','','','2',''),
    @('1244132','','','2',''),
    @('1246966','','','2',''),
    @('1264251','','','2',''),
    @('1267594','','','2',''),
    @('1272876','','','2',''),
    @('1278602','','','2',''),
    @('1290730','','','2',''),
    @('1290822','','','2',''),
    @('1299630','','','2',''),
    @('1307555','','','4',''),
    @('1310158','','','2',''),
    @('1312519','','','2',''),
    @('1317092','','','2',''),
    @('1320832','','','5',''),
    @('1328533','','','2',''),
    @('133151','','','2',''),
    @('1345526','','','3',''),
    @('1345877','','','3',''),
    @('1350143','','','2',''),
    @('1350786','','','5',''),
    @('1351207','','','2',''),
    @('1353406','','','5',''),
    @('1354859','','','3',''),
    @('1365394','','','4',''),
    @('1368285','','','2',''),
    @('136955','','','2',''),
    @('1385645','','','2',''),
    @('1387361','','','2',''),
    @('1389021','','','2',''),
    @('1402198','','','2',''),
    @('1406416','','','2',''),
    @('1406717','','','2',''),
    @('1412709','','','3',''),
    @('1417133','','','3',''),
    @('1422330','','','2',''),
    @('1438876','','','2',''),
    @('1448523','','','2',''),
    @('1462298','','','2',''),
    @('1466040','','','2',''),
    @('1470930','','','4',''),
    @('1478443','','','2',''),
    @('1478663','','','2',''),
    @('1489678','','','2',''),
    @('1491886','','','3',''),
    @('1496544','','','2',''),
    @('1497862','','','2',''),
    @('1503694','','','4',''),
    @('260606','','','2',''),
    @('307261','','','2',''),
    @('380275','','','2',''),
    @('421466','','','6',''),
    @('446429','','','5',''),
    @('469358','','','6',''),
    @('486442','','','2',''),
    @('492043','','','2',''),
    @('493192','','','2',''),
    @('510573','','','2',''),
    @('533246','','','3',''),
    @('545355','','','7',''),
    @('581579','','','2',''),
    @('621809','','','5',''),
    @('684433','','','2',''),
    @('689351','','','2',''),
    @('706455','','','2',''),
    @('722492','','','3',''),
    @('761621','','','2',''),
    @('764151','','','2',''),
    @('771041','','','3',''),
    @('785952','','','2',''),
    @('804990','','','2',''),
    @('866983','','','3',''),
    @('868202','','','2',''),
    @('872066','','','2',''),
    @('875307','','','2',''),
    @('876141','','','5',''),
    @('879581','','','2',''),
    @('880642','','','2',''),
    @('881907','','','2',''),
    @('905169','','','2',''),
    @('909642','','','2',''),
    @('937471','','','2',''),
    @('940836','','','2',''),
    @('941353','','','2',''),
    @('944163','','','2',''),
    @('950557','','','2',''),
    @('950712','','','2',''),
    @('950732','','','3',''),
    @('976980','','','2',''),
    @('982044','','','2',''),
    @('990259','','','2',''),
    @('991741','','','2',''),
    @('999592','','','3','')
)

$rowCount = $data.Count
$startRow = 2
$endRow = $startRow + $rowCount - 1

# Columns A,B,C,D,F,G hold text (SKU, titles, slugs, urls) in the source data,
# even when the text looks numeric (e.g. SKU "10051"). Force text format on
# those columns so Excel does not silently coerce them to numbers. Column E
# (Review Count) stays numeric.
$ws.Range("A$startRow`:D$endRow").NumberFormat = "@"
$ws.Range("F$startRow`:G$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $sku = $data[$i][0]
    $title = $data[$i][1]
    $slug = $data[$i][2]
    $reviewCount = [int]$data[$i][3]
    $url = $data[$i][4]
    $feefo = "https://www.feefo.com/en-US/reviews/notonthehighstreet-com/products/*?sku=$sku&displayFeedbackType=PRODUCT&timeFrame=ALL"

    $ws.Cells.Item($r, 1).Value = $sku
    $ws.Cells.Item($r, 2).Value = $title
    $ws.Cells.Item($r, 3).Value = ""
    $ws.Cells.Item($r, 4).Value = $slug
    $ws.Cells.Item($r, 5).Value = $reviewCount
    $ws.Cells.Item($r, 6).Value = $url
    $ws.Cells.Item($r, 7).Value = $feefo
}
